$wb = $excel.ActiveWorkbook

# --- Sheet 1: Julie ---
$ws1 = $wb.Worksheets.Item("Julie")
$ws1.Range("A1").Value = "Tag Name"
$ws1.Range("B1").Value = "from"
$ws1.Range("C1").Value = "subject"
$ws1.Range("D1").Value = "has"
$ws1.Range("A2").Value = "Attach"
$ws1.Range("D2").Value = "attachment"

# --- Sheet 2: Michael ---
$ws2 = $wb.Worksheets.Item("Michael")
$ws2.Range("A1").Value = "Tag Name"
$ws2.Range("B1").Value = "from"
$ws2.Range("C1").Value = "subject"
$ws2.Range("D1").Value = "has"
$ws2.Range("A2").Value = "Important"
$ws2.Range("C2").Value = "important"
$ws2.Range("A3").Value = "Trump"
$ws2.Range("C3").Value = "trump"

# --- Sheet 3: Brian ---
$ws3 = $wb.Worksheets.Item("Brian")
$ws3.Range("A1").Value = "Tag Name"
$ws3.Range("B1").Value = "from"
$ws3.Range("C1").Value = "subject"
$ws3.Range("D1").Value = "has"
$ws3.Range("A2").Value = "Foley"
$ws3.Range("A3").Value = "Target"

$ws3.Hyperlinks.Add($ws3.Range("B2"), "mailto:foleyb25@gmail.com", "", "", "foleyb25@gmail.com") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("B3"), "mailto:orders@oe.target.com", "", "", "orders@oe.target.com") | Out-Null
